$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 166, shifting existing rows 166-198 down to 167-199.
$ws.Rows("166:166").Insert()

# Populate the newly inserted row 166 with the new daily price record.
$ws.Range("A166").Value = 7
$ws.Range("B166").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C166").Value = "Ñuble"
$ws.Range("D166").Value = 45211
$ws.Range("E166").Value = 16
$ws.Range("F166").Value = "Fruta"
$ws.Range("G166").Value = 100108
$ws.Range("H166").Value = "Tropicales y subtropicales"
$ws.Range("I166").Value = 100108002
$ws.Range("J166").Value = "Mango"
$ws.Range("K166").Value = "Sin especificar"
$ws.Range("L166").Value = "Primera"
$ws.Range("M166").Value = 80
$ws.Range("N166").Value = 10000
$ws.Range("O166").Value = 10000
$ws.Range("P166").Value = 10000
$ws.Range("Q166").Value = '$/bandeja 4 kilos'
$ws.Range("R166").Value = "Brasil"
$ws.Range("S166").Value = 2500
$ws.Range("T166").Value = 4
